$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New headers in J1:L1 (percentage columns)
$ws.Range("J1").Value = "Percentuale deceduti su Positivi"
$ws.Range("K1").Value = "Percentuale Guariti su Positivi"
$ws.Range("L1").Value = "Percentuale Terpaia Intensiva su Positivi"

# Match the header style used by I1 (bold, centered, bordered)
$src = $ws.Range("I1")
$dst = $ws.Range("J1:L1")
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.LineStyle = $src.Borders.LineStyle

# Fill in missing H values (Totale Positivi precursor column) for rows 2-5
$ws.Cells.Item(2, 8).Value = 129
$ws.Cells.Item(3, 8).Value = 229
$ws.Cells.Item(4, 8).Value = 322
$ws.Cells.Item(5, 8).Value = 400

# New percentage values (J: deceduti/positivi, K: guariti/positivi, L: terapia intensiva/positivi)
$values = @(
    @(1.55, 0.78, 20.16),
    @(2.62, 0.44, 11.79),
    @(3.11, 0.31, 10.87),
    @(3, 0.75, 9),
    @(2.62, 6.92, 8.619999999999999),
    @(2.36, 5.18, 7.21),
    @(2.57, 4.43, 9.31),
    @(2.01, 4.9, 8.26),
    @(2.55, 7.32, 8.15),
    @(3.16, 6.39, 9.15),
    @(3.46, 8.93, 9.550000000000001)
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 10).Value = $v[0]
    $ws.Cells.Item($row, 11).Value = $v[1]
    $ws.Cells.Item($row, 12).Value = $v[2]
    $row++
}
